# Rename "Sheet2" to "cleaned"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "cleaned"

# Update the selection on the "cleaned" sheet to C25
$ws.Activate()
$ws.Range("C25").Select()
